# Update column F ("want-to-go count") values to match the new scrape snapshot,
# on worksheets "展览" (index 1) and "全部类型" (index 4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "展览"
$ws1.Range("F4").Value = 3764
$ws1.Range("F5").Value = 3764
$ws1.Range("F6").Value = 296
$ws1.Range("F7").Value = 5306
$ws1.Range("F8").Value = 595
$ws1.Range("F9").Value = 431
$ws1.Range("F11").Value = 1066
$ws1.Range("F13").Value = 144
$ws1.Range("F14").Value = 53
$ws1.Range("F15").Value = 733
$ws1.Range("F16").Value = 365
$ws1.Range("F17").Value = 48
$ws1.Range("F19").Value = 181
$ws1.Range("F20").Value = 5
$ws1.Range("F21").Value = 372
$ws1.Range("F22").Value = 6075
$ws1.Range("F23").Value = 6075
$ws1.Range("F25").Value = 45
$ws1.Range("F27").Value = 7030
$ws1.Range("F28").Value = 25
$ws1.Range("F29").Value = 24
$ws1.Range("F30").Value = 3262
$ws1.Range("F31").Value = 374
$ws1.Range("F32").Value = 754
$ws1.Range("F33").Value = 4462
$ws1.Range("F35").Value = 137
$ws1.Range("F37").Value = 1168
$ws1.Range("F38").Value = 105
$ws1.Range("F39").Value = 29
$ws1.Range("F41").Value = 928
$ws1.Range("F42").Value = 1152
$ws1.Range("F43").Value = 2063

$ws4 = $wb.Worksheets.Item(4)   # "全部类型"
$ws4.Range("F7").Value = 3764
$ws4.Range("F8").Value = 3764
$ws4.Range("F9").Value = 296
$ws4.Range("F10").Value = 5306
$ws4.Range("F11").Value = 595
$ws4.Range("F12").Value = 431
$ws4.Range("F14").Value = 1066
$ws4.Range("F16").Value = 144
$ws4.Range("F17").Value = 53
$ws4.Range("F18").Value = 733
$ws4.Range("F19").Value = 365
$ws4.Range("F20").Value = 48
$ws4.Range("F23").Value = 181
$ws4.Range("F24").Value = 5
$ws4.Range("F25").Value = 372
$ws4.Range("F26").Value = 6075
$ws4.Range("F28").Value = 45
$ws4.Range("F30").Value = 7030
$ws4.Range("F31").Value = 25
$ws4.Range("F32").Value = 24
$ws4.Range("F33").Value = 3262
$ws4.Range("F34").Value = 374
$ws4.Range("F35").Value = 754
$ws4.Range("F36").Value = 4462
$ws4.Range("F39").Value = 137
$ws4.Range("F41").Value = 1168
$ws4.Range("F42").Value = 105
$ws4.Range("F43").Value = 29
$ws4.Range("F45").Value = 928
$ws4.Range("F46").Value = 1152
$ws4.Range("F48").Value = 2063
